$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (styles) from the last existing data row (2699) down to the new rows (2700-2729)
$ws.Range("A2699:H2699").Copy()
$ws.Range("A2700:H2729").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(2700,1).Value = 44180
$ws.Cells.Item(2700,2).Value = "0-10 years"
$ws.Cells.Item(2700,3).Value = 24077
$ws.Cells.Item(2700,4).Value = [double]"5.09162040708432E-02"
$ws.Cells.Item(2700,5).Value = 515
$ws.Cells.Item(2700,6).Value = [double]"6.24166767664526E-02"
$ws.Cells.Item(2700,7).Value = 4
$ws.Cells.Item(2700,8).Value = 0

$ws.Cells.Item(2701,1).Value = 44180
$ws.Cells.Item(2701,2).Value = "11-20 years"
$ws.Cells.Item(2701,3).Value = 61932
$ws.Cells.Item(2701,4).Value = 0.130969072164948
$ws.Cells.Item(2701,5).Value = 1121
$ws.Cells.Item(2701,6).Value = 0.135862319718822
$ws.Cells.Item(2701,7).Value = 3
$ws.Cells.Item(2701,8).Value = 0

$ws.Cells.Item(2702,1).Value = 44180
$ws.Cells.Item(2702,2).Value = "21-30 years"
$ws.Cells.Item(2702,3).Value = 88470
$ws.Cells.Item(2702,4).Value = 0.187089611419508
$ws.Cells.Item(2702,5).Value = 1455
$ws.Cells.Item(2702,6).Value = 0.176342261544055
$ws.Cells.Item(2702,7).Value = 31
$ws.Cells.Item(2702,8).Value = 0

$ws.Cells.Item(2703,1).Value = 44180
$ws.Cells.Item(2703,2).Value = "31-40 years"
$ws.Cells.Item(2703,3).Value = 73976
$ws.Cells.Item(2703,4).Value = 0.156438805181073
$ws.Cells.Item(2703,5).Value = 1189
$ws.Cells.Item(2703,6).Value = 0.144103745000606
$ws.Cells.Item(2703,7).Value = 64
$ws.Cells.Item(2703,8).Value = 0

$ws.Cells.Item(2704,1).Value = 44180
$ws.Cells.Item(2704,2).Value = "41-50 years"
$ws.Cells.Item(2704,3).Value = 70487
$ws.Cells.Item(2704,4).Value = 0.14906053396775
$ws.Cells.Item(2704,5).Value = 1143
$ws.Cells.Item(2704,6).Value = 0.13852866319234
$ws.Cells.Item(2704,7).Value = 179
$ws.Cells.Item(2704,8).Value = 3

$ws.Cells.Item(2705,1).Value = 44180
$ws.Cells.Item(2705,2).Value = "51-60 years"
$ws.Cells.Item(2705,3).Value = 64902
$ws.Cells.Item(2705,4).Value = 0.137249801744647
$ws.Cells.Item(2705,5).Value = 1208
$ws.Cells.Item(2705,6).Value = 0.146406496182281
$ws.Cells.Item(2705,7).Value = 497
$ws.Cells.Item(2705,8).Value = 10

$ws.Cells.Item(2706,1).Value = 44180
$ws.Cells.Item(2706,2).Value = "61-70 years"
$ws.Cells.Item(2706,3).Value = 46118
$ws.Cells.Item(2706,4).Value = [double]"9.75268305577584E-02"
$ws.Cells.Item(2706,5).Value = 799
$ws.Cells.Item(2706,6).Value = [double]"9.68367470609623E-02"
$ws.Cells.Item(2706,7).Value = 1042
$ws.Cells.Item(2706,8).Value = 6

$ws.Cells.Item(2707,1).Value = 44180
$ws.Cells.Item(2707,2).Value = "71-80 years"
$ws.Cells.Item(2707,3).Value = 27645
$ws.Cells.Item(2707,4).Value = [double]"5.84615384615385E-02"
$ws.Cells.Item(2707,5).Value = 526
$ws.Cells.Item(2707,6).Value = [double]"6.37498485032117E-02"
$ws.Cells.Item(2707,7).Value = 1703
$ws.Cells.Item(2707,8).Value = 22

$ws.Cells.Item(2708,1).Value = 44180
$ws.Cells.Item(2708,2).Value = "81+ years"
$ws.Cells.Item(2708,3).Value = 14701
$ws.Cells.Item(2708,4).Value = [double]"3.10885540576262E-02"
$ws.Cells.Item(2708,5).Value = 295
$ws.Cells.Item(2708,6).Value = [double]"3.57532420312689E-02"
$ws.Cells.Item(2708,7).Value = 2092
$ws.Cells.Item(2708,8).Value = 33

$ws.Cells.Item(2709,1).Value = 44180
$ws.Cells.Item(2709,2).Value = "Pending"
$ws.Cells.Item(2709,3).Value = 567
$ws.Cells.Item(2709,4).Value = [double]"1.19904837430611E-03"
$ws.Cells.Item(2709,5).Value = 0
$ws.Cells.Item(2709,6).Value = 0
$ws.Cells.Item(2709,7).Value = 0
$ws.Cells.Item(2709,8).Value = 0

$ws.Cells.Item(2710,1).Value = 44181
$ws.Cells.Item(2710,2).Value = "0-10 years"
$ws.Cells.Item(2710,3).Value = 24695
$ws.Cells.Item(2710,4).Value = [double]"5.09927005792044E-02"
$ws.Cells.Item(2710,5).Value = 618
$ws.Cells.Item(2710,6).Value = [double]"5.41630148992112E-02"
$ws.Cells.Item(2710,7).Value = 4
$ws.Cells.Item(2710,8).Value = 0

$ws.Cells.Item(2711,1).Value = 44181
$ws.Cells.Item(2711,2).Value = "11-20 years"
$ws.Cells.Item(2711,3).Value = 63328
$ws.Cells.Item(2711,4).Value = 0.130765974581083
$ws.Cells.Item(2711,5).Value = 1396
$ws.Cells.Item(2711,6).Value = 0.122348816827344
$ws.Cells.Item(2711,7).Value = 3
$ws.Cells.Item(2711,8).Value = 0

$ws.Cells.Item(2712,1).Value = 44181
$ws.Cells.Item(2712,2).Value = "21-30 years"
$ws.Cells.Item(2712,3).Value = 90258
$ws.Cells.Item(2712,4).Value = 0.186373726214935
$ws.Cells.Item(2712,5).Value = 1788
$ws.Cells.Item(2712,6).Value = 0.156704645048203
$ws.Cells.Item(2712,7).Value = 31
$ws.Cells.Item(2712,8).Value = 0

$ws.Cells.Item(2713,1).Value = 44181
$ws.Cells.Item(2713,2).Value = "31-40 years"
$ws.Cells.Item(2713,3).Value = 75682
$ws.Cells.Item(2713,4).Value = 0.15627574671939
$ws.Cells.Item(2713,5).Value = 1706
$ws.Cells.Item(2713,6).Value = 0.149517966695881
$ws.Cells.Item(2713,7).Value = 64
$ws.Cells.Item(2713,8).Value = 0

$ws.Cells.Item(2714,1).Value = 44181
$ws.Cells.Item(2714,2).Value = "41-50 years"
$ws.Cells.Item(2714,3).Value = 72271
$ws.Cells.Item(2714,4).Value = 0.149232373499076
$ws.Cells.Item(2714,5).Value = 1784
$ws.Cells.Item(2714,6).Value = 0.15635407537248
$ws.Cells.Item(2714,7).Value = 179
$ws.Cells.Item(2714,8).Value = 0

$ws.Cells.Item(2715,1).Value = 44181
$ws.Cells.Item(2715,2).Value = "51-60 years"
$ws.Cells.Item(2715,3).Value = 66557
$ws.Cells.Item(2715,4).Value = 0.137433536037664
$ws.Cells.Item(2715,5).Value = 1655
$ws.Cells.Item(2715,6).Value = 0.145048203330412
$ws.Cells.Item(2715,7).Value = 501
$ws.Cells.Item(2715,8).Value = 4

$ws.Cells.Item(2716,1).Value = 44181
$ws.Cells.Item(2716,2).Value = "61-70 years"
$ws.Cells.Item(2716,3).Value = 47309
$ws.Cells.Item(2716,4).Value = [double]"9.76883446730747E-02"
$ws.Cells.Item(2716,5).Value = 1191
$ws.Cells.Item(2716,6).Value = 0.104382120946538
$ws.Cells.Item(2716,7).Value = 1049
$ws.Cells.Item(2716,8).Value = 7

$ws.Cells.Item(2717,1).Value = 44181
$ws.Cells.Item(2717,2).Value = "71-80 years"
$ws.Cells.Item(2717,3).Value = 28440
$ws.Cells.Item(2717,4).Value = [double]"5.87257503329651E-02"
$ws.Cells.Item(2717,5).Value = 795
$ws.Cells.Item(2717,6).Value = [double]"6.96757230499562E-02"
$ws.Cells.Item(2717,7).Value = 1720
$ws.Cells.Item(2717,8).Value = 17

$ws.Cells.Item(2718,1).Value = 44181
$ws.Cells.Item(2718,2).Value = "81+ years"
$ws.Cells.Item(2718,3).Value = 15184
$ws.Cells.Item(2718,4).Value = [double]"3.13534385743932E-02"
$ws.Cells.Item(2718,5).Value = 483
$ws.Cells.Item(2718,6).Value = [double]"4.23312883435583E-02"
$ws.Cells.Item(2718,7).Value = 2117
$ws.Cells.Item(2718,8).Value = 25

$ws.Cells.Item(2719,1).Value = 44181
$ws.Cells.Item(2719,2).Value = "Pending"
$ws.Cells.Item(2719,3).Value = 561
$ws.Cells.Item(2719,4).Value = [double]"1.15840878821355E-03"
$ws.Cells.Item(2719,5).Value = -6
$ws.Cells.Item(2719,6).Value = [double]"-5.25854513584575E-04"
$ws.Cells.Item(2719,7).Value = 0
$ws.Cells.Item(2719,8).Value = 0

$ws.Cells.Item(2720,1).Value = 44182
$ws.Cells.Item(2720,2).Value = "0-10 years"
$ws.Cells.Item(2720,3).Value = 25176
$ws.Cells.Item(2720,4).Value = [double]"5.10431238975731E-02"
$ws.Cells.Item(2720,5).Value = 481
$ws.Cells.Item(2720,6).Value = [double]"5.37730575740637E-02"
$ws.Cells.Item(2720,7).Value = 4
$ws.Cells.Item(2720,8).Value = 0

$ws.Cells.Item(2721,1).Value = 44182
$ws.Cells.Item(2721,2).Value = "11-20 years"
$ws.Cells.Item(2721,3).Value = 64447
$ws.Cells.Item(2721,4).Value = 0.13066317944975
$ws.Cells.Item(2721,5).Value = 1119
$ws.Cells.Item(2721,6).Value = 0.125097820011179
$ws.Cells.Item(2721,7).Value = 4
$ws.Cells.Item(2721,8).Value = 1

$ws.Cells.Item(2722,1).Value = 44182
$ws.Cells.Item(2722,2).Value = "21-30 years"
$ws.Cells.Item(2722,3).Value = 91850
$ws.Cells.Item(2722,4).Value = 0.186221438274233
$ws.Cells.Item(2722,5).Value = 1592
$ws.Cells.Item(2722,6).Value = 0.177976523197317
$ws.Cells.Item(2722,7).Value = 31
$ws.Cells.Item(2722,8).Value = 0

$ws.Cells.Item(2723,1).Value = 44182
$ws.Cells.Item(2723,2).Value = "31-40 years"
$ws.Cells.Item(2723,3).Value = 77121
$ws.Cells.Item(2723,4).Value = 0.156359102244389
$ws.Cells.Item(2723,5).Value = 1439
$ws.Cells.Item(2723,6).Value = 0.160871995528228
$ws.Cells.Item(2723,7).Value = 66
$ws.Cells.Item(2723,8).Value = 2

$ws.Cells.Item(2724,1).Value = 44182
$ws.Cells.Item(2724,2).Value = "41-50 years"
$ws.Cells.Item(2724,3).Value = 73709
$ws.Cells.Item(2724,4).Value = 0.149441437057762
$ws.Cells.Item(2724,5).Value = 1438
$ws.Cells.Item(2724,6).Value = 0.160760201229737
$ws.Cells.Item(2724,7).Value = 183
$ws.Cells.Item(2724,8).Value = 4

$ws.Cells.Item(2725,1).Value = 44182
$ws.Cells.Item(2725,2).Value = "51-60 years"
$ws.Cells.Item(2725,3).Value = 67839
$ws.Cells.Item(2725,4).Value = 0.137540295602457
$ws.Cells.Item(2725,5).Value = 1282
$ws.Cells.Item(2725,6).Value = 0.143320290665176
$ws.Cells.Item(2725,7).Value = 507
$ws.Cells.Item(2725,8).Value = 6

$ws.Cells.Item(2726,1).Value = 44182
$ws.Cells.Item(2726,2).Value = "61-70 years"
$ws.Cells.Item(2726,3).Value = 48210
$ws.Cells.Item(2726,4).Value = [double]"9.77434462623928E-02"
$ws.Cells.Item(2726,5).Value = 901
$ws.Cells.Item(2726,6).Value = 0.10072666294019
$ws.Cells.Item(2726,7).Value = 1079
$ws.Cells.Item(2726,8).Value = 30

$ws.Cells.Item(2727,1).Value = 44182
$ws.Cells.Item(2727,2).Value = "71-80 years"
$ws.Cells.Item(2727,3).Value = 28941
$ws.Cells.Item(2727,4).Value = [double]"5.86764795328751E-02"
$ws.Cells.Item(2727,5).Value = 501
$ws.Cells.Item(2727,6).Value = [double]"5.60089435438793E-02"
$ws.Cells.Item(2727,7).Value = 1769
$ws.Cells.Item(2727,8).Value = 49

$ws.Cells.Item(2728,1).Value = 44182
$ws.Cells.Item(2728,2).Value = "81+ years"
$ws.Cells.Item(2728,3).Value = 15388
$ws.Cells.Item(2728,4).Value = [double]"3.11984266974839E-02"
$ws.Cells.Item(2728,5).Value = 204
$ws.Cells.Item(2728,6).Value = [double]"2.28060368921185E-02"
$ws.Cells.Item(2728,7).Value = 2202
$ws.Cells.Item(2728,8).Value = 85

$ws.Cells.Item(2729,1).Value = 44182
$ws.Cells.Item(2729,2).Value = "Pending"
$ws.Cells.Item(2729,3).Value = 549
$ws.Cells.Item(2729,4).Value = [double]"1.11307098108388E-03"
$ws.Cells.Item(2729,5).Value = -12
$ws.Cells.Item(2729,6).Value = [double]"-1.34153158188932E-03"
$ws.Cells.Item(2729,7).Value = 0
$ws.Cells.Item(2729,8).Value = 0

# Update the defined name range to cover the new rows
$wb.Names.Item("ALL_AGE_FINAL").RefersTo = "='ALL_AGE_FINAL'!`$A`$1:`$H`$2729"
